# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# - Updates "VALOR MORA" and "Cant. Periodos" summary cells
# - Replaces the worker ledger (rows 16-26) with the updated data set:
#     * ANDRES DIONICIO YEPEZ BARRIOS (CC 1007968914) periods 2409-2501
#     * KELINETH DEL CARMEN PAREDES FONTALVO (CC 45547507) periods 2502-2508
#   (one new trailing period, 2508, for Kelineth)
# - Inserts a new row so the bottom-bordered "last row" style follows the
#   new last data row, and the signature/footer rows shift down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Summary header updates -------------------------------------------------
$ws.Range("E11").Value = 658580
$ws.Range("F13").Value = 12

# ---- Make room for the extra ledger row -------------------------------------
# Before: last data row is 26, footer (signature) rows are 31-32.
# After:  last data row is 27, footer (signature) rows are 32-33.
$ws.Rows(27).Insert()

# The freshly inserted row 27 picked up a blended style; give it the exact
# "last row" (bottom border) formatting that row 26 currently has...
$ws.Range("B26:J26").Copy()
$ws.Range("B27:J27").PasteSpecial(-4122)

# ...and restore row 26 back to the regular (non-bottom-border) row style,
# copied from row 25, since it is no longer the final ledger row.
$ws.Range("B25:J25").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# ---- Rewrite the worker ledger rows 16-27 -----------------------------------
$tipoDoc = "CC"

$andres = "ANDRES DIONICIO YEPEZ BARRIOS"
$andresDoc = 1007968914
$andresPeriods = @("2409", "2410", "2411", "2412", "2501")
$andresSalario = 52000
$andresBasico = 1423500

$kelineth = "KELINETH DEL CARMEN PAREDES FONTALVO"
$kelinethDoc = 45547507
$kelinethPeriods = @("2502", "2503", "2504", "2505", "2506", "2507", "2508")
$kelinethSalario = 56940
$kelinethBasico = 1423500

$row = 16
foreach ($periodo in $andresPeriods) {
    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $andresDoc
    $ws.Cells.Item($row, 4).Value = $andres
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $andresSalario
    $ws.Cells.Item($row, 7).Value = $andresBasico
    $row = $row + 1
}

foreach ($periodo in $kelinethPeriods) {
    $ws.Cells.Item($row, 2).Value = $tipoDoc
    $ws.Cells.Item($row, 3).Value = $kelinethDoc
    $ws.Cells.Item($row, 4).Value = $kelineth
    $ws.Cells.Item($row, 5).Value = $periodo
    $ws.Cells.Item($row, 6).Value = $kelinethSalario
    $ws.Cells.Item($row, 7).Value = $kelinethBasico
    $row = $row + 1
}
